# Append one new sensor-log data row (columns A-I) to each of the 4
# worksheets, continuing directly after the last existing row. Each new
# row is seeded by copying the row immediately above it (so formatting -
# e.g. the datetime style on column A, and the text/number typing on the
# other columns - is inherited exactly, without introducing any new style
# entries), then the per-column values are overwritten on top.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1 - ROW50-FE-LIFTER: new row 98 ----------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A97:I97").Copy($ws1.Range("A98:I98"))
$ws1.Cells.Item(98, 1).Value = 45773.3124137963
$ws1.Cells.Item(98, 2).Value = "0x01,0x90"
$ws1.Cells.Item(98, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item(98, 4).Value = "0x01,0x3a"
$ws1.Cells.Item(98, 5).Value = "0xe"
$ws1.Cells.Item(98, 6).Value = 400
$ws1.Cells.Item(98, 7).Value = $ws1.Cells.Item(2, 7).Value2
$ws1.Cells.Item(98, 8).Value = 314
$ws1.Cells.Item(98, 9).Value = 14

# ---- Sheet 2 - ROW50-MID-LIFTER: new row 100 ---------------------------
# (the ID_DEC column, G, is stored as text throughout this sheet; row 99's
# G value already equals the needed constant, so leave it untouched by the
# copy instead of re-typing it as a numeric-looking string)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A99:I99").Copy($ws2.Range("A100:I100"))
$ws2.Cells.Item(100, 1).Value = 45773.27555555556
$ws2.Cells.Item(100, 2).Value = "0x01,0x90 "
$ws2.Cells.Item(100, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item(100, 4).Value = "0x01,0x3e"
$ws2.Cells.Item(100, 5).Value = "0x19"
$ws2.Cells.Item(100, 6).Value = 400
$ws2.Cells.Item(100, 8).Value = 318
$ws2.Cells.Item(100, 9).Value = 25

# ---- Sheet 3 - ROW11-FE-LIFTER: new row 98 -----------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A97:I97").Copy($ws3.Range("A98:I98"))
$ws3.Cells.Item(98, 1).Value = 45773.34379238426
$ws3.Cells.Item(98, 2).Value = "0x01,0x90"
$ws3.Cells.Item(98, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item(98, 4).Value = "0x01,0x3a"
$ws3.Cells.Item(98, 5).Value = "0x14"
$ws3.Cells.Item(98, 6).Value = 400
$ws3.Cells.Item(98, 7).Value = $ws3.Cells.Item(2, 7).Value2
$ws3.Cells.Item(98, 8).Value = 314
$ws3.Cells.Item(98, 9).Value = 20

# ---- Sheet 4 - ROW11-MID-LIFTER: new row 98 ----------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A97:I97").Copy($ws4.Range("A98:I98"))
$ws4.Cells.Item(98, 1).Value = 45773.46272855324
$ws4.Cells.Item(98, 2).Value = "0x01,0x90"
$ws4.Cells.Item(98, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item(98, 4).Value = "0x01,0x3e"
$ws4.Cells.Item(98, 5).Value = "0x19"
$ws4.Cells.Item(98, 6).Value = 400
$ws4.Cells.Item(98, 7).Value = $ws4.Cells.Item(2, 7).Value2
$ws4.Cells.Item(98, 8).Value = 318
$ws4.Cells.Item(98, 9).Value = 25
